$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" summary sheet: insert a new data row for "2022-Q3" right after the
#    header, pushing the existing quarters down by one row.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Range("A2:D2").EntireRow.Insert()

$summaryData = @(
    @(0, "2022-Q3", 13, 0.78),
    @(1, "2022-Q2", 9, 0.78),
    @(2, "2022-Q1", 13, 1.57),
    @(3, "2021-Q4", 7, 0.92),
    @(4, "2021-Q3", 6, 0.46),
    @(5, "2021-Q2", 4, 0.49),
    @(6, "2021-Q1", 5, 0.52),
    @(7, "2020-Q4", 2, 0.67)
)

for ($i = 0; $i -lt $summaryData.Length; $i++) {
    $r = $i + 2
    $row = $summaryData[$i]
    $summary.Cells.Item($r, 1).Value = $row[0]
    $summary.Cells.Item($r, 2).Value = $row[1]
    $summary.Cells.Item($r, 3).Value = $row[2]
    $summary.Cells.Item($r, 4).Value = $row[3]
}

# ---------------------------------------------------------------------------
# 2) New "2022-Q3" detail sheet, holding the same fund-level columns as the
#    other quarterly sheets. Duplicate the "2022-Q2" sheet (placing the copy
#    right before it) so the new sheet inherits identical column/style
#    formatting, then rename it and overwrite its contents.
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q2")
$template.Copy($template)
$newSheet = $wb.Worksheets.Item("2022-Q2 (2)")
$newSheet.Name = "2022-Q3"

$dataQ3 = @(
    @(0, "'159745", "'国泰中证全指建筑材料ETF", "'7.92", "'99.14", "'3.33", "'0.2637", 8),
    @(1, "'004856", "'广发中证全指建筑材料指数A", "'7.66", "'93.74", "'3.14", "'0.2405", 8),
    @(2, "'004857", "'广发中证全指建筑材料指数C", "'6.12", "'93.74", "'3.14", "'0.1922", 8),
    @(3, "'516750", "'富国中证全指建筑材料ETF", "'0.82", "'98.46", "'3.25", "'0.0266", 8),
    @(4, "'012419", "'天弘国证建材指数C", "'0.63", "'94.93", "'2.80", "'0.0176", 10),
    @(5, "'009658", "'汇丰晋信中小盘低波动策略股票A", "'0.85", "'90.14", "'1.95", "'0.0166", 5),
    @(6, "'003242", "'创金合信量化发现灵活配置混合C", "'0.40", "'92.08", "'1.50", "'0.0060", 8),
    @(7, "'159787", "'易方达中证全指建筑材料ETF", "'0.17", "'94.24", "'3.20", "'0.0054", 8),
    @(8, "'003241", "'创金合信量化发现灵活配置混合A", "'0.32", "'92.08", "'1.50", "'0.0048", 8),
    @(9, "'012405", "'天弘国证建材指数A", "'0.13", "'94.93", "'2.80", "'0.0036", 10),
    @(10, "'164811", "'工银瑞信中证京津冀协同发展主题指数（LOF）A", "'0.12", "'93.09", "'2.87", "'0.0034", 10),
    @(11, "'164825", "'工银瑞信中证京津冀协同发展主题指数（LOF）C", "'0.03", "'93.09", "'2.87", "'0.0009", 10),
    @(12, "'009775", "'汇丰晋信中小盘低波动策略股票C", "'0.04", "'90.14", "'1.95", "'0.0008", 5)
)

# Extend the A-column / row styling (inherited from the template's data rows)
# down to the extra rows the new sheet needs beyond the template's row count.
$newSheet.Range("A2:H2").Copy()
$newSheet.Range("A11:H14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $dataQ3.Length; $i++) {
    $r = $i + 2
    $row = $dataQ3[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $newSheet.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}

$newSheet.Range("A1").Select()
Write-Output "2022-Q3 sheet created and 总计 updated"
